$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-43 from
# serial date 45841 (2025-07-03) to 45842 (2025-07-04).
$ws.Range("C2:C43").Value = 45842
